$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 60000
$ws.Range("J87").Value = 60000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62496

$ws.Range("H90").Value = 60000
$ws.Range("J90").Value = 60000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -192480

$ws.Range("H98").Value = 5264.25
$ws.Range("I98").Value = 5058.7827
$ws.Range("K98").Value = 5058.7827
$ws.Range("M98").Value = -3560.7827

$ws.Range("H116").Value = 19239844
$ws.Range("I116").Value = 27783666
$ws.Range("K116").Value = 27783666
$ws.Range("M116").Value = -27780224

$ws.Range("H122").Value = 5264.25
$ws.Range("I122").Value = 5058.7827
$ws.Range("K122").Value = 15176.3481
$ws.Range("M122").Value = -12726.3481

$ws.Range("H132").Value = 1868.7317
$ws.Range("I132").Value = 1412.5883
$ws.Range("J132").Value = 4084.2856
$ws.Range("K132").Value = 4237.7649
$ws.Range("L132").Value = 12252.8568
$ws.Range("M132").Value = -1707.7649
$ws.Range("N132").Value = -17312.8568

$ws.Range("H137").Value = 2024.64
$ws.Range("I137").Value = 2386.4614
$ws.Range("K137").Value = 7159.3842
$ws.Range("M137").Value = -4609.3842

$ws.Range("H138").Value = 6613.3213
$ws.Range("J138").Value = 12106.929
$ws.Range("L138").Value = 36320.787
$ws.Range("N138").Value = -46600.787

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2457.7097
$ws.Range("I2").Value = 1757.7333
$ws.Range("J2").Value = 3113.9375
$ws.Range("K2").Value = 1757.7333
$ws.Range("L2").Value = 3113.9375
$ws.Range("M2").Value = -1644.7333
$ws.Range("N2").Value = -3339.9375

$ws.Range("H32").Value = 1691527.4
$ws.Range("I32").Value = 1787925.9
$ws.Range("K32").Value = 1787925.9
$ws.Range("M32").Value = -1787638.9

$ws.Range("H61").Value = 6362.243
$ws.Range("I61").Value = 2810.5173
$ws.Range("J61").Value = 19237.25
$ws.Range("K61").Value = 2810.5173
$ws.Range("L61").Value = 19237.25
$ws.Range("M61").Value = -2598.5173
$ws.Range("N61").Value = -19661.25

$ws.Range("H88").Value = 3890.182
$ws.Range("J88").Value = 4799
$ws.Range("L88").Value = 4799
$ws.Range("N88").Value = -5611

$ws.Range("H91").Value = 3890.182
$ws.Range("J91").Value = 4799
$ws.Range("L91").Value = 4799
$ws.Range("N91").Value = -7607

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H116").Value = 2457.7097
$ws.Range("I116").Value = 1757.7333
$ws.Range("J116").Value = 3113.9375
$ws.Range("K116").Value = 1757.7333
$ws.Range("L116").Value = 3113.9375
$ws.Range("M116").Value = 536.2666999999999
$ws.Range("N116").Value = -7701.9375

$ws.Range("H132").Value = 11763.315
$ws.Range("I132").Value = 14658.429
$ws.Range("K132").Value = 43975.287
$ws.Range("M132").Value = -41445.287

$ws.Range("H136").Value = 6362.243
$ws.Range("I136").Value = 2810.5173
$ws.Range("J136").Value = 19237.25
$ws.Range("K136").Value = 8431.5519
$ws.Range("L136").Value = 57711.75
$ws.Range("M136").Value = -5881.5519
$ws.Range("N136").Value = -62811.75

$ws.Range("H138").Value = 79807.836
$ws.Range("I138").Value = 79424
$ws.Range("J138").Value = 79884.60000000001
$ws.Range("K138").Value = 79424
$ws.Range("L138").Value = 79884.60000000001
$ws.Range("M138").Value = -74284
$ws.Range("N138").Value = -90164.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2457.7097
$ws.Range("I3").Value = 1757.7333
$ws.Range("J3").Value = 3113.9375
$ws.Range("K3").Value = 1757.7333
$ws.Range("L3").Value = 3113.9375
$ws.Range("M3").Value = -1643.7333
$ws.Range("N3").Value = -3341.9375

$ws.Range("H105").Value = 1757.8572
$ws.Range("I105").Value = 1238.5264
$ws.Range("K105").Value = 1238.5264
$ws.Range("M105").Value = 508.4736

$ws.Range("H134").Value = 6853.143
$ws.Range("I134").Value = 2956.1177
$ws.Range("J134").Value = 10533.667
$ws.Range("K134").Value = 8868.3531
$ws.Range("L134").Value = 31601.001
$ws.Range("M134").Value = -6333.3531
$ws.Range("N134").Value = -36671.001

$ws.Range("H139").Value = 29999.75
$ws.Range("J139").Value = 29999.75
$ws.Range("L139").Value = 29999.75
$ws.Range("N139").Value = -40279.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 8045.476
$ws.Range("I132").Value = 3639.4285
$ws.Range("J132").Value = 10248.5
$ws.Range("K132").Value = 10918.2855
$ws.Range("L132").Value = 30745.5
$ws.Range("M132").Value = -8388.2855
$ws.Range("N132").Value = -35805.5

$ws.Range("H134").Value = 8668.625
$ws.Range("I134").Value = 4186.091
$ws.Range("K134").Value = 12558.273
$ws.Range("M134").Value = -10023.273

$ws.Range("H141").Value = 60295.547
$ws.Range("J141").Value = 60295.547
$ws.Range("L141").Value = 60295.547
$ws.Range("N141").Value = -70655.54699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 39218890
$ws.Range("J75").Value = 27781374
$ws.Range("L75").Value = 83344122
$ws.Range("N75").Value = -83346118

$ws.Range("H78").Value = 39218890
$ws.Range("J78").Value = 27781374
$ws.Range("L78").Value = 250032366
$ws.Range("N78").Value = -250042350

$ws.Range("H92").Value = 9617122
$ws.Range("J92").Value = 9617122
$ws.Range("L92").Value = 28851366
$ws.Range("N92").Value = -28853862

$ws.Range("H98").Value = 1179.7693
$ws.Range("I98").Value = 893.1667
$ws.Range("J98").Value = 1425.4286
$ws.Range("K98").Value = 2679.5001
$ws.Range("L98").Value = 4276.2858
$ws.Range("M98").Value = -1181.5001
$ws.Range("N98").Value = -7272.2858

$ws.Range("H132").Value = 9190.393
$ws.Range("I132").Value = 3771.75
$ws.Range("J132").Value = 16415.25
$ws.Range("K132").Value = 33945.75
$ws.Range("L132").Value = 147737.25
$ws.Range("M132").Value = -31415.75
$ws.Range("N132").Value = -152797.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4497.5
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 7495
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 7495
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -9491

$ws.Range("H83").Value = 4497.5
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 7495
$ws.Range("K83").Value = 7500
$ws.Range("L83").Value = 37475
$ws.Range("M83").Value = -2508
$ws.Range("N83").Value = -47459

$ws.Range("H97").Value = 6785.7144
$ws.Range("I97").Value = 6166.6665
$ws.Range("J97").Value = 7250
$ws.Range("K97").Value = 6166.6665
$ws.Range("L97").Value = 7250
$ws.Range("M97").Value = -5670.6665
$ws.Range("N97").Value = -8242

$ws.Range("H102").Value = 10733.125
$ws.Range("I102").Value = 10310.833
$ws.Range("K102").Value = 10310.833
$ws.Range("M102").Value = -8688.833000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2803.2
$ws.Range("I22").Value = 1499.5
$ws.Range("J22").Value = 3129.125
$ws.Range("K22").Value = 1499.5
$ws.Range("L22").Value = 3129.125
$ws.Range("M22").Value = -1204.5
$ws.Range("N22").Value = -3719.125

$ws.Range("H27").Value = 2803.2
$ws.Range("I27").Value = 1499.5
$ws.Range("J27").Value = 3129.125
$ws.Range("K27").Value = 1499.5
$ws.Range("L27").Value = 3129.125
$ws.Range("M27").Value = -1392.5
$ws.Range("N27").Value = -3343.125

$ws.Range("H61").Value = 3227982.8
$ws.Range("I61").Value = 4546515.5
$ws.Range("J61").Value = 4902.6665
$ws.Range("K61").Value = 4546515.5
$ws.Range("L61").Value = 4902.6665
$ws.Range("M61").Value = -4546313.5
$ws.Range("N61").Value = -5306.6665

$ws.Range("H82").Value = 1086759
$ws.Range("I82").Value = 1762997.5
$ws.Range("K82").Value = 1762997.5
$ws.Range("M82").Value = -1762636.5

$ws.Range("H85").Value = 1086759
$ws.Range("I85").Value = 1762997.5
$ws.Range("K85").Value = 1762997.5
$ws.Range("M85").Value = -1761749.5

$ws.Range("H113").Value = 3227982.8
$ws.Range("I113").Value = 4546515.5
$ws.Range("J113").Value = 4902.6665
$ws.Range("K113").Value = 4546515.5
$ws.Range("L113").Value = 4902.6665
$ws.Range("M113").Value = -4546313.5
$ws.Range("N113").Value = -9242.666499999999
